$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.003343645993482
$ws.Range("C2").Value = 0.4297108598773036
$ws.Range("D2").Value = 0.5638575298801101
$ws.Range("E2").Value = 0.1855828576089635
$ws.Range("G2").Value = 0.002696976618000235
$ws.Range("J2").Value = 0.06684150712439418
$ws.Range("K2").Value = 2.353406400852435
$ws.Range("M2").Value = 0.8176558072228417
$ws.Range("N2").Value = 6.107543588151685
$ws.Range("B3").Value = 1.955813302402845
$ws.Range("C3").Value = 0.420724349350138
$ws.Range("D3").Value = 0.5602427652440269
$ws.Range("E3").Value = 0.1848138450510497
$ws.Range("G3").Value = 0.002702845023073719
$ws.Range("J3").Value = 0.06681583488865073
$ws.Range("K3").Value = 2.29921239858686
$ws.Range("M3").Value = 0.804826314824453
$ws.Range("N3").Value = 6.028648744157408
$ws.Range("B4").Value = 1.927814312756595
$ws.Range("C4").Value = 0.4154621514965413
$ws.Range("D4").Value = 0.5582932737252122
$ws.Range("E4").Value = 0.1844258204972569
$ws.Range("G4").Value = 0.002706635370396706
$ws.Range("J4").Value = 0.06682657586843277
$ws.Range("K4").Value = 2.267330988947634
$ws.Range("M4").Value = 0.7974030618743839
$ws.Range("N4").Value = 5.980628484069626
$ws.Range("B5").Value = 1.916701440399038
$ws.Range("C5").Value = 0.4133817592172591
$ws.Range("D5").Value = 0.5575666526629703
$ws.Range("E5").Value = 0.1842888399981106
$ws.Range("G5").Value = 0.002708227189086517
$ws.Range("J5").Value = 0.06683761572539737
$ws.Range("K5").Value = 2.254688332927401
$ws.Range("M5").Value = 0.7944919203135825
$ws.Range("N5").Value = 5.96116463182139
$ws.Range("B6").Value = 1.914874062554333
$ws.Range("C6").Value = 0.4130401692226542
$ws.Range("D6").Value = 0.5574500912996001
$ws.Range("E6").Value = 0.1842673709508844
$ws.Range("G6").Value = 0.002708494366505105
$ws.Range("J6").Value = 0.0668398512579067
$ws.Range("K6").Value = 2.252610089472682
$ws.Range("M6").Value = 0.7940153979734745
$ws.Range("N6").Value = 5.957938968551019
$ws.Range("B7").Value = 1.927663239669897
$ws.Range("C7").Value = 0.4154338358281962
$ws.Range("D7").Value = 0.5582831997880788
$ws.Range("E7").Value = 0.1844238875473607
$ws.Range("G7").Value = 0.002706656646802432
$ws.Range("J7").Value = 0.0668266977804528
$ws.Range("K7").Value = 2.267159073084429
$ws.Range("M7").Value = 0.7973633404338472
$ws.Range("N7").Value = 5.980365565089414
$ws.Range("B8").Value = 1.986708647187243
$ws.Range("C8").Value = 0.4265591051818092
$ws.Range("D8").Value = 0.5625550722033523
$ws.Range("E8").Value = 0.1853002243758191
$ws.Range("G8").Value = 0.002698961308174886
$ws.Range("J8").Value = 0.06682715315174548
$ws.Range("K8").Value = 2.334430224011186
$ws.Range("M8").Value = 0.8131377625520742
$ws.Range("N8").Value = 6.08025241922897
$ws.Range("B9").Value = 2.111951920610466
$ws.Range("C9").Value = 0.4504173006019414
$ws.Range("D9").Value = 0.573079248891645
$ws.Range("E9").Value = 0.187687600024006
$ws.Range("G9").Value = 0.002685347834063514
$ws.Range("J9").Value = 0.06703850742852779
$ws.Range("K9").Value = 2.477475668405418
$ws.Range("M9").Value = 0.8476907582222921
$ws.Range("N9").Value = 6.279539017560751
$ws.Range("B10").Value = 2.209821354979965
$ws.Range("C10").Value = 0.4692124182340933
$ws.Range("D10").Value = 0.5821287765295153
$ws.Range("E10").Value = 0.1898515098050488
$ws.Range("G10").Value = 0.002676235704824864
$ws.Range("J10").Value = 0.06732244454240188
$ws.Range("K10").Value = 2.589462677704375
$ws.Range("M10").Value = 0.8753101572821649
$ws.Range("N10").Value = 6.42813740137268
$ws.Range("B11").Value = 2.255636122203384
$ws.Range("C11").Value = 0.4780428052752654
$ws.Range("D11").Value = 0.5865337327015254
$ws.Range("E11").Value = 0.1909254505886224
$ws.Range("G11").Value = 0.002672281246678882
$ws.Range("J11").Value = 0.06747964125212391
$ws.Range("K11").Value = 2.641929737119938
$ws.Range("M11").Value = 0.8883660834804417
$ws.Range("N11").Value = 6.49623613297166
$ws.Range("B12").Value = 2.273172429197245
$ws.Range("C12").Value = 0.4814273296323393
$ws.Range("D12").Value = 0.5882433697153431
$ws.Range("E12").Value = 0.1913450391980263
$ws.Range("G12").Value = 0.002670811041439104
$ws.Range("J12").Value = 0.06754320463883445
$ws.Range("K12").Value = 2.662018514745114
$ws.Range("M12").Value = 0.8933811729324788
$ws.Range("N12").Value = 6.52209693001538
$ws.Range("B13").Value = 2.269387323522324
$ws.Range("C13").Value = 0.4806965985965235
$ws.Range("D13").Value = 0.5878733177185893
$ws.Range("E13").Value = 0.1912540985766817
$ws.Range("G13").Value = 0.002671126466597923
$ws.Range("J13").Value = 0.06752933551578977
$ws.Range("K13").Value = 2.657682198908503
$ws.Range("M13").Value = 0.8922979155426773
$ws.Range("N13").Value = 6.516524064355622
$ws.Range("B14").Value = 2.25707508471146
$ws.Range("C14").Value = 0.4783204355678379
$ws.Range("D14").Value = 0.586673551614922
$ws.Range("E14").Value = 0.1909597114674177
$ws.Range("G14").Value = 0.002672159746533309
$ws.Range("J14").Value = 0.06748478972231098
$ws.Range("K14").Value = 2.643578021964231
$ws.Range("M14").Value = 0.8887772508457559
$ws.Range("N14").Value = 6.498362240050881
$ws.Range("B15").Value = 2.249557905918891
$ws.Range("C15").Value = 0.4768702706011254
$ws.Range("D15").Value = 0.5859440783067384
$ws.Range("E15").Value = 0.1907810731379662
$ws.Range("G15").Value = 0.002672796205966577
$ws.Range("J15").Value = 0.06745802996347905
$ws.Range("K15").Value = 2.63496758691798
$ws.Range("M15").Value = 0.8866300121395057
$ws.Range("N15").Value = 6.48724718106962
$ws.Range("B16").Value = 2.2068533861883
$ws.Range("C16").Value = 0.468641001785727
$ws.Range("D16").Value = 0.5818467134820082
$ws.Range("E16").Value = 0.1897831301783519
$ws.Range("G16").Value = 0.002676497961580855
$ws.Range("J16").Value = 0.06731273590101949
$ws.Range("K16").Value = 2.586064622282322
$ws.Range("M16").Value = 0.8744668491293695
$ws.Range("N16").Value = 6.423697173060589
$ws.Range("B17").Value = 2.180987714010826
$ws.Range("C17").Value = 0.4636646585019264
$ws.Range("D17").Value = 0.5794070331697583
$ws.Range("E17").Value = 0.1891938863420286
$ws.Range("G17").Value = 0.002678817596729843
$ws.Range("J17").Value = 0.06723078615465283
$ws.Range("K17").Value = 2.556455558777486
$ws.Range("M17").Value = 0.8671313546344237
$ws.Range("N17").Value = 6.384840434221076
$ws.Range("B18").Value = 2.166232134378333
$ws.Range("C18").Value = 0.4608287693590967
$ws.Range("D18").Value = 0.5780309203621243
$ws.Range("E18").Value = 0.1888633974427059
$ws.Range("G18").Value = 0.002680169748372765
$ws.Range("J18").Value = 0.0671862890136552
$ws.Range("K18").Value = 2.539568522748311
$ws.Range("M18").Value = 0.8629584318629142
$ws.Range("N18").Value = 6.362538104096416
$ws.Range("B19").Value = 2.161257008196287
$ws.Range("C19").Value = 0.459873105548354
$ws.Range("D19").Value = 0.577569647646726
$ws.Range("E19").Value = 0.1887529462784165
$ws.Range("G19").Value = 0.002680630652822216
$ws.Range("J19").Value = 0.06717167599915541
$ws.Range("K19").Value = 2.533875430318687
$ws.Range("M19").Value = 0.8615534855176605
$ws.Range("N19").Value = 6.354994957745163
$ws.Range("B20").Value = 2.183728557553025
$ws.Range("C20").Value = 0.4641916669970101
$ws.Range("D20").Value = 0.5796639329234949
$ws.Range("E20").Value = 0.1892557398105517
$ws.Range("G20").Value = 0.002678568810284489
$ws.Range("J20").Value = 0.06723923676521792
$ws.Range("K20").Value = 2.55959265057777
$ws.Range("M20").Value = 0.8679074401004314
$ws.Range("N20").Value = 6.388971922733163
$ws.Range("B21").Value = 2.260686395697974
$ws.Range("C21").Value = 0.4790172664075669
$ws.Range("D21").Value = 0.587024822657952
$ws.Range("E21").Value = 0.1910458294558026
$ws.Range("G21").Value = 0.00267185550838988
$ws.Range("J21").Value = 0.06749776430987353
$ws.Range("K21").Value = 2.647714759957182
$ws.Range("M21").Value = 0.8898094228320375
$ws.Range("N21").Value = 6.503694808695798
$ws.Range("B22").Value = 2.312074825481488
$ws.Range("C22").Value = 0.4889437407385628
$ws.Range("D22").Value = 0.5920779846434812
$ws.Range("E22").Value = 0.1922910244261189
$ws.Range("G22").Value = 0.002667626810415767
$ws.Range("J22").Value = 0.06769025704792853
$ws.Range("K22").Value = 2.706594463588374
$ws.Range("M22").Value = 0.904538208673344
$ws.Range("N22").Value = 6.579100625038393
$ws.Range("B23").Value = 2.284547510725076
$ws.Range("C23").Value = 0.4836239905190212
$ws.Range("D23").Value = 0.5893587964122275
$ws.Range("E23").Value = 0.191619542832079
$ws.Range("G23").Value = 0.002669869265641153
$ws.Range("J23").Value = 0.06758536516316482
$ws.Range("K23").Value = 2.675050992460683
$ws.Range("M23").Value = 0.8966391149981234
$ws.Range("N23").Value = 6.538815545571424
$ws.Range("B24").Value = 2.182489064269532
$ws.Range("C24").Value = 0.463953328372213
$ws.Range("D24").Value = 0.5795477060376868
$ws.Range("E24").Value = 0.1892277500842106
$ws.Range("G24").Value = 0.002678681228737325
$ws.Range("J24").Value = 0.067235408093552
$ws.Range("K24").Value = 2.558173949099341
$ws.Range("M24").Value = 0.867556433747481
$ws.Range("N24").Value = 6.387103961856582
$ws.Range("B25").Value = 2.077048950332653
$ws.Range("C25").Value = 0.4437422244332083
$ws.Range("D25").Value = 0.5700014611144297
$ws.Range("E25").Value = 0.1869699274756655
$ws.Range("G25").Value = 0.002688873625100598
$ws.Range("J25").Value = 0.06695876412478796
$ws.Range("K25").Value = 2.437575817537635
$ws.Range("M25").Value = 0.8379530403461715
$ws.Range("N25").Value = 6.225251032835388
